# Remove the "Ver no Jupiter..." block (an empty paragraph, the text
# paragraph itself, another empty paragraph, and the following empty
# page-break paragraph) that used to follow the "LOQ4031: Química Geral I
# (Requisito fraco)" requirement line.

$d = $word.ActiveDocument

$marker = "Ver no Jupiter Salvar em pdf Salvar em docx"

# Locate the anchor paragraph ("LOQ4031: ...") by scanning the document's
# paragraphs for its distinctive text.
$anchorPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOQ4031*") {
        $anchorPara = $p
        break
    }
}

if ($anchorPara -eq $null) {
    throw "Anchor paragraph not found"
}

# The next four paragraphs are: an empty paragraph, the "Ver no Jupiter..."
# paragraph, another empty paragraph, and the empty page-break paragraph.
$p1 = $anchorPara.Next()
$p2 = $p1.Next()
$p3 = $p2.Next()
$p4 = $p3.Next()

if ($p2.Range.Text -notlike "*Ver no Jupiter*") {
    throw "Unexpected paragraph layout while locating block to remove"
}

$deleteRange = $d.Range($p1.Range.Start, $p4.Range.End)
$deleteRange.Delete()

Write-Output "Removed trailing 'Ver no Jupiter' block"
